# "Add files via upload" - add Problem 14 / Problem 15 data sheets, plus a
# blank Sheet2, and refresh a handful of sheet-view states (zoom/selection)
# left over from the author's last interactive session.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Problem 14: repurpose the old blank "Sheet4" tab, filling in an
# item / monthly-stock table. Column A (item names) is entered before
# the header row, matching how the shared-string table ended up ordered.
# ------------------------------------------------------------------
$ws14 = $wb.Worksheets.Item("Sheet4")
$ws14.Name = "Problem 14"

$ws14.Range("A2").Value = "Pen"
$ws14.Range("A3").Value = "Notebook"

# ------------------------------------------------------------------
# Problem 15: brand-new sheet with a student score table.
# ------------------------------------------------------------------
$ws15 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws14)
$ws15.Name = "Problem 15"

$ws15.Range("A1").Value = "student_id"
$ws15.Range("B1").Value = "name"
$ws15.Range("C1").Value = "subject "
$ws15.Range("D1").Value = "score"

$ws15.Range("B2").Value = "Alice"
$ws15.Range("B3").Value = "Bob"
$ws15.Range("B4").Value = "Carol"
$ws15.Range("B5").Value = "David"
$ws15.Range("B6").Value = "Emily"

$ws15.Range("C2").Value = "Math"
$ws15.Range("C3").Value = "Science"
$ws15.Range("C4").Value = "English"
$ws15.Range("C5").Value = "Math"
$ws15.Range("C6").Value = "Science"

$ws15.Range("A2").Value = 1
$ws15.Range("A3").Value = 2
$ws15.Range("A4").Value = 3
$ws15.Range("A5").Value = 4
$ws15.Range("A6").Value = 5

$ws15.Range("D2").Value = 91
$ws15.Range("D3").Value = 88
$ws15.Range("D4").Value = 92
$ws15.Range("D5").Value = 85
$ws15.Range("D6").Value = 95

# ------------------------------------------------------------------
# Back to Problem 14: header row + the numeric stock columns.
# ------------------------------------------------------------------
$ws14.Range("A1").Value = "item"
$ws14.Range("B1").Value = "jan_stock"
$ws14.Range("C1").Value = "feb_stock"
$ws14.Range("D1").Value = "mar_stock"
$ws14.Range("E1").Value = "apr_stock"

$ws14.Range("B2").Value = 150
$ws14.Range("C2").Value = 120
$ws14.Range("D2").Value = 175
$ws14.Range("E2").Value = 200

$ws14.Range("B3").Value = 250
$ws14.Range("C3").Value = 280
$ws14.Range("D3").Value = 300
$ws14.Range("E3").Value = 270

# ------------------------------------------------------------------
# A new, still-empty "Sheet2" tab added after Problem 15.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws15)
$ws2.Name = "Sheet2"

# ------------------------------------------------------------------
# Restore the per-sheet view state (zoom level + last selection) that
# was left behind on a few tabs.
# ------------------------------------------------------------------
$wsDf1 = $wb.Worksheets.Item("df1 Problem 12")
$wsDf1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 202
$wsDf1.Range("C8").Select() | Out-Null

$wsDf2 = $wb.Worksheets.Item("df2 Problem 12")
$wsDf2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 172
$wsDf2.Range("C12").Select() | Out-Null

$ws13 = $wb.Worksheets.Item("Problem 13")
$ws13.Activate() | Out-Null
$ws13.Range("F9").Select() | Out-Null

$ws14.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 219
$ws14.Range("E12").Select() | Out-Null

$ws15.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 211
$ws15.Range("E8").Select() | Out-Null
